# Add Russian ("ru") translation column (F) to the translations sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('F1').Value  = 'ru'
$ws.Range('F2').Value  = 'Вид кадрирования'
$ws.Range('F3').Value  = 'Выберите вид кадрирования'
$ws.Range('F4').Value  = 'Обрезать холст'
$ws.Range('F5').Value  = 'Создать маску'
$ws.Range('F6').Value  = 'Фоновый слой'
$ws.Range('F7').Value  = 'Фоновая заливка'
$ws.Range('F8').Value  = ' - сетка'
$ws.Range('F9').Value  = ' - размер'
$ws.Range('F10').Value = ' - показ'
$ws.Range('F11').Value = ' - обрезка'
$ws.Range('F12').Value = 'Golden Crop от SzopeN'
$ws.Range('F13').Value = 'Маска кадрирвания'
$ws.Range('F14').Value = 'Правила разделения'
$ws.Range('F15').Value = 'Линии на %1%%'
$ws.Range('F16').Value = 'Золотой треугольник вверх'
$ws.Range('F17').Value = 'Золотой треугольник вниз'
$ws.Range('F18').Value = 'Метод диагоналей'
$ws.Range('F19').Value = 'Откройте файл, в котором вы хотели бы запустить скрипт.'
$ws.Range('F20').Value = 'Обнаружено увеличение холста'
$ws.Range('F21').Value = 'Что следует предпринять?'
$ws.Range('F22').Value = 'Увеличить размеры'
$ws.Range('F23').Value = 'Обрезать без увеличения'
$ws.Range('F24').Value = 'Вернуться к кадрированию'
$ws.Range('F25').Value = 'Метод построения композиции'
$ws.Range('F26').Value = 'Выберите тип направляющих линий'
$ws.Range('F27').Value = 'Золотое сечение'
$ws.Range('F28').Value = 'Правило третей'
$ws.Range('F29').Value = 'Золотая спираль внизу-слева'
$ws.Range('F30').Value = 'Золотая спираль вверху-слева'
$ws.Range('F31').Value = 'Золотая спираль вверху-справа'
$ws.Range('F32').Value = 'Золотая спираль внизу-справа'
$ws.Range('F33').Value = 'Выбрать все'
$ws.Range('F34').Value = 'Убрать все'
$ws.Range('F35').Value = 'OK'
$ws.Range('F36').Value = 'Отмена'
$ws.Range('F37').Value = 'Все золотые спирали'
$ws.Range('F38').Value = 'Основные правила'
$ws.Range('F39').Value = 'Толщина линий'
$ws.Range('F40').Value = 'Толщина линий (‰ меньшей стороны): '

# Match the new column's width (closest value the width-quantization of this
# host allows to the authored 53.42578125 characters).
$ws.Columns.Item(6).ColumnWidth = 52.65

# The author's selection ends up on the newly added header cell.
[void]$ws.Range('F1').Select()
